$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Remove the extra duplicate header cells C1:F1 on the header row
$ws.Range("C1:F1").ClearContents()

# Row 8: "Model" -> "production_function"
$ws.Range("A8").Value = "production_function"

# Insert a new row for the L_curve parameter right after row 8
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# Remove the old "Deletion" row (now shifted down to row 17)
$ws.Rows("17:17").Delete()

# Make optimization_parameters the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("C1:G7").Select()
